# Update the Angpt1-Itgb1 NATMI LR-pair sheet with the new TPM-derived values.
#
# The underlying NATMI script was re-run against new TPM numbers:
#  - rows whose "Sending cluster" was ECs are dropped entirely
#  - the remaining FAPs/MuSCs rows get refreshed receptor/edge statistics

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old "ECs"-as-sender rows (original rows 2-4); this shifts the
# FAPs/MuSCs rows up so the table runs from row 2 to row 7 again.
$ws.Range("A2:A4").EntireRow.Delete()

$data = @(
    @("FAPs",  "Angpt1", "Itgb1", "ECs",   3, 1, 12.95459633333333, 38.863789,          0.8906505749177925, 0.8906505749177924, 3, 1, 77.08952333333333,  231.26857,          0.2403816673726824, 0.2403816673726824, 998.6636563124144,  8987.97290681173,  0.2140960702451772,  0.2140960702451771),
    @("FAPs",  "Angpt1", "Itgb1", "FAPs",  3, 1, 12.95459633333333, 38.863789,          0.8906505749177925, 0.8906505749177924, 3, 1, 101.5800373333333,  304.740112,         0.3167483425780597, 0.3167483425780597, 1315.928379178263,  11843.35541260437, 0.2821120934214068,  0.2821120934214067),
    @("FAPs",  "Angpt1", "Itgb1", "MuSCs", 3, 1, 12.95459633333333, 38.863789,          0.8906505749177925, 0.8906505749177924, 3, 1, 142.0267893333333,  426.080368,         0.4428699900492579, 0.4428699900492579, 1839.899724332706,  16559.09751899435, 0.3944424112512086,  0.3944424112512086),
    @("MuSCs", "Angpt1", "Itgb1", "ECs",   3, 1, 1.590497666666667, 4.771493,           0.1093494250822076, 0.1093494250822076, 3, 1, 77.08952333333333,  231.26857,          0.2403816673726824, 0.2403816673726824, 122.6107069861122,  1103.49636287501,  0.02628559712750527, 0.02628559712750527),
    @("MuSCs", "Angpt1", "Itgb1", "FAPs",  3, 1, 1.590497666666667, 4.771493,           0.1093494250822076, 0.1093494250822076, 3, 1, 101.5800373333333,  304.740112,         0.3167483425780597, 0.3167483425780597, 161.5628123585796,  1454.065311227216, 0.03463624915665296, 0.03463624915665296),
    @("MuSCs", "Angpt1", "Itgb1", "MuSCs", 3, 1, 1.590497666666667, 4.771493,           0.1093494250822076, 0.1093494250822076, 3, 1, 142.0267893333333,  426.080368,         0.4428699900492579, 0.4428699900492579, 225.8932770388249,  2033.039493349424, 0.04842757879804934, 0.04842757879804934)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $values = $data[$i]
    for ($c = 1; $c -le $values.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}
